$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting (bold, border, centered)
# from an existing header cell, then set its text.
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (H2:H17)
$saveValues = @(1, 1, 0, 0, 0, 0, 0, 0, 1, 0, 1, 0, 1, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
